$wb = $excel.ActiveWorkbook

$renames = @{
    "3_500_500_outliers_upregulate" = "3_500_500_fp_upregulated"
    "3_500_500_outliers_downregula" = "3_500_500_fp_downregulated"
    "3_750_250_outliers_upregulate" = "3_750_250_fp_upregulated"
    "3_750_250_outliers_downregula" = "3_750_250_fp_downregulated"
    "3_1000_0_outliers_upregulated" = "3_1000_0_fp_upregulated"
    "3_1000_0_outliers_downregulat" = "3_1000_0_fp_downregulated"
    "6_500_500_outliers_upregulate" = "6_500_500_fp_upregulated"
    "6_500_500_outliers_downregula" = "6_500_500_fp_downregulated"
    "6_750_250_outliers_upregulate" = "6_750_250_fp_upregulated"
    "6_750_250_outliers_downregula" = "6_750_250_fp_downregulated"
    "6_1000_0_outliers_upregulated" = "6_1000_0_fp_upregulated"
    "6_1000_0_outliers_downregulat" = "6_1000_0_fp_downregulated"
    "9_500_500_outliers_upregulate" = "9_500_500_fp_upregulated"
    "9_500_500_outliers_downregula" = "9_500_500_fp_downregulated"
    "9_750_250_outliers_upregulate" = "9_750_250_fp_upregulated"
    "9_750_250_outliers_downregula" = "9_750_250_fp_downregulated"
    "9_1000_0_outliers_upregulated" = "9_1000_0_fp_upregulated"
    "9_1000_0_outliers_downregulat" = "9_1000_0_fp_downregulated"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renames.ContainsKey($oldName)) {
        $ws.Name = $renames[$oldName]
    }
}
